# Update 2p0. Convention change to support multi-axle vehicles
#
# Rename the "Scalable" sedan variant sheet to "HambaLG" to reflect the new
# multi-axle naming convention, update the matching label cell that mirrors
# the variant name, and restore the user's last selection in the frozen
# pane (the bottom-right pane was last parked on H4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab.
$ws.Name = "Sedan_HambaLG"

# The variant-name label cell (H3) mirrors the sheet/tab name - update it too.
$ws.Range("H3").Value = "Sedan_HambaLG"

# Restore the last active selection in the bottom-right (frozen) pane.
$ws.Range("H4").Select() | Out-Null
